$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" (row 2) ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.300711743772242
$ws1.Range("C2").Value = 0.06443914081145585
$ws1.Range("E2").Value = 0.1208053691275168
$ws1.Range("F2").Value = 0.2542372881355932
$ws1.Range("G2").Value = 0.6273458445040214
$ws1.Range("H2").Value = 0.7626404494382023
$ws1.Range("J2").Value = 392
$ws1.Range("K2").Value = 142

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# Row 2 - label "0"
$ws2.Range("B2").Value = 0.993006993006993
$ws2.Range("C2").Value = 0.2659176029962547
$ws2.Range("D2").Value = 0.4194977843426883

# Row 3 - label "1"
$ws2.Range("B3").Value = 0.06443914081145585
$ws2.Range("D3").Value = 0.1208053691275168

# Row 4 - label "accuracy"
$ws2.Range("B4").Value = 0.300711743772242
$ws2.Range("C4").Value = 0.300711743772242
$ws2.Range("D4").Value = 0.300711743772242
$ws2.Range("E4").Value = 0.300711743772242

# Row 5 - label "macro avg"
$ws2.Range("B5").Value = 0.5287230669092244
$ws2.Range("C5").Value = 0.6151016586409845
$ws2.Range("D5").Value = 0.2701515767351025

# Row 6 - label "weighted avg"
$ws2.Range("B6").Value = 0.9467438259936922
$ws2.Range("C6").Value = 0.300711743772242
$ws2.Range("D6").Value = 0.404616311698516

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 142
$ws3.Range("C2").Value = 392
